# Update countries & provincias Spain
# - Refresh the COVID-19 case counters for the countries whose ranking
#   shifted in today's feed (column order: B=Casos totales, C=Nuevos casos,
#   D=Casos activos, E=Recuperados, F=Casos criticos, G=Muertes hoy,
#   H=Muertes).
# - A few countries leap-frogged their neighbour in the ranking, so their
#   row's country name (column A) needs to be swapped too.
# - Bump the "Datos actualizados" timestamp footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => new B..H values (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes)
$rowData = @{
    4   = @(6515506, 1275, 3797553, 2523878, 0, 43, 194075)
    5   = @(4382518, 15082, 3406270, 902220, 0, 105, 74028)
    13  = @(500034, 0, 382490, 107087, 0, 52, 10457)
    19  = @(323012, 775, 298966, 19881, 0, 28, 4165)
    24  = @(255041, 85, 230600, 15032, 0, 0, 9409)
    27  = @(143030, 2551, 64703, 75348, 0, 45, 2979)
    41  = @(85880, 0, 0, 0, 0, 4, 5842)
    44  = @(77688, 1140, 0, 0, 0, 2, 6246)
    47  = @(73402, 194, 71999, 677, 0, 5, 726)
    51  = @(61541, 646, 43284, 16408, 0, 3, 1849)
    60  = @(45306, 469, 38100, 5187, 0, 1, 2019)
    63  = @(44781, 500, 42162, 2255, 0, 6, 364)
    64  = @(44613, 87, 40336, 3216, 0, 1, 1061)
    67  = @(37732, 175, 35119, 2058, 0, 3, 555)
    69  = @(32078, 84, 30780, 570, 0, 1, 728)
    77  = @(22258, 297, 15208, 6375, 0, 6, 675)
    80  = @(20462, 879, 2329, 17809, 0, 10, 324)
    81  = @(19848, 0, 18448, 985, 0, 0, 415)
    83  = @(18607, 251, 15990, 1989, 0, 0, 628)
    86  = @(15293, 67, 12754, 1905, 0, 3, 634)
    89  = @(13112, 160, 11839, 973, 0, 3, 300)
    94  = @(10324, 32, 9523, 541, 0, 0, 260)
    146 = @(2153, 3, 2067, 76, 0, 0, 10)
    194 = @(116, 9, 51, 64, 0, 0, 1)
    195 = @(108, 1, 105, 2, 0, 0, 1)
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, $i + 2).Value = $vals[$i]
    }
}

# Countries that overtook their neighbour in the ranking: swap the labels.
$labelSwaps = @{
    63  = "Uzbekistan"
    64  = "Kirguistan"
    80  = "Libia"
    81  = "Camerun"
    194 = "Curazao"
    195 = "Liechtenstein"
}

foreach ($row in $labelSwaps.Keys) {
    $ws.Cells.Item($row, 1).Value = $labelSwaps[$row]
}

# Footer timestamp.
$ws.Range("A1").Value = "Datos actualizados a 9 de Septiembre de 2020 a las 15:45"
